# Daily attendance processing - reverse the order of names/emails listed
# in the "Recorded By" column (G) wherever multiple comma-separated
# values are present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val.ToString().Contains(",")) {
        $parts = $val.ToString().Split(",")

        $trimmed = @()
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $trimmed += $parts[$i].Trim()
        }

        # [array]::Reverse() is unreliable in this runtime, so reverse manually
        $n = $trimmed.Length
        $reversed = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
